# Refresh the cryptos price list (rows 2-51 of Sheet1) with the latest
# scraped snapshot: updated Price / Volume(1h) figures, and - where the
# underlying source re-ranked a coin - the Coin name + Link swapped between
# two adjacent rows.
#
# Price-looking strings (e.g. "314.12", "1.003") are forced to Text via
# NumberFormat "@" before assignment so Excel doesn't reinterpret them as
# numbers, matching the original inline-string storage of this column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.359.20'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.858.20'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.12'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4641'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  +1.55%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07348'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8828'
$ws.Range('E10').Value = '  +2.79%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.990.61'
$ws.Range('E11').Value = '  +15.78%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07892'
$ws.Range('E12').Value = '  +2.69%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.89'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.397'
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.575'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '92.07'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008882'
$ws.Range('E18').Value = '  +2.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('E20').Value = '  +2.60%  '
$ws.Range('D21').Value = '27.400.64'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.127'
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '2.132.93'
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.96'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.887'
$ws.Range('E26').Value = '  +2.68%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.40'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.080'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.130'
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '116.33'
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08893'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7578'
$ws.Range('E32').Value = '  +5.31%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.025'
$ws.Range('E33').Value = '  +2.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.163'
$ws.Range('E34').Value = '  +2.94%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.494'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.612'
$ws.Range('E36').Value = '  +8.19%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.077'
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01957'
$ws.Range('E38').Value = '  +1.69%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.975'
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05231'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.114'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5166'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.343'
$ws.Range('E44').Value = '  +2.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4839'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.24'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.003'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '103.48'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.654'
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06236'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '65.69'
$ws.Range('E51').Value = '  +1.87%  '
